# Update data to the September/October version across the three sheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: "HSV Log OLS"
$ws1 = $wb.Worksheets.Item("HSV Log OLS")
$ws1.Range("A2").Value = 0.108
$ws1.Range("B2").Value = 3.027
$ws1.Range("C2").Value = 0.859

# Sheet 2: "HSV PPML"
$ws2 = $wb.Worksheets.Item("HSV PPML")
$ws2.Range("A2").Value = 0.04
$ws2.Range("B2").Value = 1.372
$ws2.Range("C2").Value = 0.861

# Sheet 3: "HSVT NLLSQ"
$ws3 = $wb.Worksheets.Item("HSVT NLLSQ")
$ws3.Range("A2").Value = -0.054
$ws3.Range("B2").Value = 0.392
$ws3.Range("C2").Value = 11983.85
$ws3.Range("D2").Value = 0.104
